$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 and add new rows 3-6 with the new product data
$ws.Range("A2").Value = 5151
$ws.Range("B2").Value = "ريد بل - 250 مل"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 1065
$ws.Range("E2").Value = "YES"

$ws.Range("A3").Value = 13928
$ws.Range("B3").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1065
$ws.Range("E3").Value = "YES"

$ws.Range("A4").Value = 13928
$ws.Range("B4").Value = "ريد بول 12 كانز - 250 مل"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 532
$ws.Range("E4").Value = "YES"

$ws.Range("A5").Value = 19977
$ws.Range("B5").Value = "ريدبل كريز و توت بري  - 250 مل"
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 1050
$ws.Range("E5").Value = "YES"

$ws.Range("A6").Value = 7630
$ws.Range("B6").Value = "فيورى جولد - 400 مل"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 205
$ws.Range("E6").Value = "YES"
